$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# Clear the "smoothness" values (value 1) in column C that are no longer needed.
# Rows 5, 15 and 30 carry a row-specific style that differs from the column's
# default style, so only their value is cleared (the cell itself stays, keeping
# its style). The other rows share column C's default style, so clearing them
# fully removes the cell from the sheet, matching how Excel drops "default"
# cells that hold no data.
$ws.Range("C5").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C30").ClearContents()

$ws.Range("C6:C13").Clear()
$ws.Range("C16:C28").Clear()
$ws.Range("C31").Clear()

# Update frozen pane top-left cell and the active selection to reflect scroll position.
$ws.Activate()
$ws.Range("C31").Select()
$excel.ActiveWindow.ScrollRow = 10

